$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readAntares")

# mcYears comment: clarify that multiple years can be given separated by ";"
$ws.Range("C11").Value = "one or more numbers (separated by ;) or synthetic, or empty / NULL / NA. Multiple years : 1;2"

# removeVirtualAreas comment is no longer needed -> fully clear cell (incl. formatting)
$ws.Range("C43").Clear()

# storageFlexibility rows: drop the (duplicated / misplaced) comment text but
# keep the existing cell formatting
$ws.Range("C44:C49").ClearContents()

# production row: comment was actually describing a disabled/enabled flag
$ws.Range("C50").Value = "0 disabled - 1 enabled"

# reassignCost / newCols rows: these reference virtual storage/flexibility areas
$ws.Range("C51").Value = "names of the virtual storage/flexibility areas"
$ws.Range("C52").Value = "names of the virtual storage/flexibility areas"

# leave the selection on the comment column of the parameters table
[void]$ws.Range("C1:C52").Select()
